$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 91830
$ws.Range("B3").Value = 79245
$ws.Range("B4").Value = 91773
$ws.Range("B5").Value = 79245
$ws.Range("B6").Value = 79245
$ws.Range("B7").Value = 91830
$ws.Range("B8").Value = 79245
$ws.Range("B9").Value = 91810
$ws.Range("B10").Value = 79245
$ws.Range("B11").Value = 92108
